$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 255, shifting existing rows 255:268 down to 256:269
$ws.Rows(255).Insert()

# Populate the newly inserted row 255 with the new record
$ws.Range("A255").Value = 2
$ws.Range("B255").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C255").Value = "Coquimbo"
$ws.Range("D255").Value = 45021
$ws.Range("E255").Value = 4
$ws.Range("F255").Value = 100112031
$ws.Range("G255").Value = "Poroto verde"
$ws.Range("H255").Value = "Magnum"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 700
$ws.Range("K255").Value = 14000
$ws.Range("L255").Value = 17000
$ws.Range("M255").Value = 15500
$ws.Range("N255").Value = "$/malla 25 kilos"
$ws.Range("O255").Value = "Provincia de Limarí"
$ws.Range("P255").Value = 620
$ws.Range("Q255").Value = 25
$ws.Range("R255").Value = "Hortaliza"
